# B6-PowerPoint.pptx edit script
#
# Per the commit diff:
#   1. ppt/presentation.xml gains embedTrueTypeFonts="1" plus a
#      <p:embeddedFontLst> entry embedding the "Tahoma" font (regular +
#      bold).  Embedding real font binaries is a PowerPoint "Save
#      Options" / file-format feature -- it has never been exposed on
#      the Presentation/Application COM object model (no
#      EmbedTrueTypeFonts property, no writable Fonts.Add, etc.), so it
#      cannot be produced from VBA/COM automation. There is nothing to
#      script for that part.
#   2. Three tables (on the slides whose graphicFrame holds a table)
#      have their <a:tableStyleId> switched from
#      {56272B1C-A215-48B7-8F59-2619C03564CA} (the deck's custom
#      "Table_0" style) to {CC33BB10-E678-441E-B173-BC4872E3414F} (a
#      built-in PowerPoint table style GUID). That part *is* reachable
#      through automation via Table.ApplyStyle(styleId).

$p = $ppt.ActivePresentation

$oldStyleId = "{56272B1C-A215-48B7-8F59-2619C03564CA}"
$newStyleId = "{CC33BB10-E678-441E-B173-BC4872E3414F}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
